# Restored from revision of admin on 09/17/2020 07:13:29 AM.TEST Author: admin. Type: SAVE.
# The only substantive content change is the "Integer min" value for the
# R40 rule row (row 10) in the "Rules" sheet, moving from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
